$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13.30571766666667
$ws.Range("H2").Value = 39.917153
$ws.Range("I2").Value = 0.007643519924167935
$ws.Range("J2").Value = 0.007643519924167933
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.752937333333333
$ws.Range("N2").Value = 11.258812
$ws.Range("O2").Value = 0.6855621274031838
$ws.Range("P2").Value = 0.6855621274031838
$ws.Range("Q2").Value = 49.93552457802622
$ws.Range("R2").Value = 449.419721202236
$ws.Range("S2").Value = 0.005240107780061191
$ws.Range("T2").Value = 0.00524010778006119

$ws.Range("G3").Value = 13.30571766666667
$ws.Range("H3").Value = 39.917153
$ws.Range("I3").Value = 0.007643519924167935
$ws.Range("J3").Value = 0.007643519924167933
$ws.Range("O3").Value = 0.2368266084628361
$ws.Range("P3").Value = 0.2368266084628362
$ws.Range("Q3").Value = 17.25016662227544
$ws.Range("R3").Value = 155.251499600479
$ws.Range("S3").Value = 0.001810188900358806
$ws.Range("T3").Value = 0.001810188900358806

$ws.Range("G4").Value = 13.30571766666667
$ws.Range("H4").Value = 39.917153
$ws.Range("I4").Value = 0.007643519924167935
$ws.Range("J4").Value = 0.007643519924167933
$ws.Range("M4").Value = 0.4248633333333334
$ws.Range("N4").Value = 1.27459
$ws.Range("O4").Value = 0.07761126413398003
$ws.Range("P4").Value = 0.07761126413398005
$ws.Range("Q4").Value = 5.653111560252222
$ws.Range("R4").Value = 50.87800404227
$ws.Range("S4").Value = 0.0005932232437479366
$ws.Range("T4").Value = 0.0005932232437479366

$ws.Range("I5").Value = 0.9448263940026712
$ws.Range("J5").Value = 0.9448263940026712
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.752937333333333
$ws.Range("N5").Value = 11.258812
$ws.Range("O5").Value = 0.6855621274031838
$ws.Range("P5").Value = 0.6855621274031838
$ws.Range("Q5").Value = 6172.601378392335
$ws.Range("R5").Value = 55553.41240553102
$ws.Range("S5").Value = 0.64773719269915
$ws.Range("T5").Value = 0.64773719269915

$ws.Range("I6").Value = 0.9448263940026712
$ws.Range("J6").Value = 0.9448263940026712
$ws.Range("O6").Value = 0.2368266084628361
$ws.Range("P6").Value = 0.2368266084628362
$ws.Range("S6").Value = 0.2237600304778239
$ws.Range("T6").Value = 0.223760030477824

$ws.Range("I7").Value = 0.9448263940026712
$ws.Range("J7").Value = 0.9448263940026712
$ws.Range("M7").Value = 0.4248633333333334
$ws.Range("N7").Value = 1.27459
$ws.Range("O7").Value = 0.07761126413398003
$ws.Range("P7").Value = 0.07761126413398005
$ws.Range("Q7").Value = 698.7891787237488
$ws.Range("R7").Value = 6289.10260851374
$ws.Range("S7").Value = 0.0733291708256972
$ws.Range("T7").Value = 0.07332917082569722

$ws.Range("G8").Value = 82.73961633333333
$ws.Range("H8").Value = 248.218849
$ws.Range("I8").Value = 0.04753008607316088
$ws.Range("J8").Value = 0.04753008607316087
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 3.752937333333333
$ws.Range("N8").Value = 11.258812
$ws.Range("O8").Value = 0.6855621274031838
$ws.Range("P8").Value = 0.6855621274031838
$ws.Range("Q8").Value = 310.5165950830431
$ws.Range("R8").Value = 2794.649355747388
$ws.Range("S8").Value = 0.03258482692397261
$ws.Range("T8").Value = 0.0325848269239726

$ws.Range("G9").Value = 82.73961633333333
$ws.Range("H9").Value = 248.218849
$ws.Range("I9").Value = 0.04753008607316088
$ws.Range("J9").Value = 0.04753008607316087
$ws.Range("O9").Value = 0.2368266084628361
$ws.Range("P9").Value = 0.2368266084628362
$ws.Range("Q9").Value = 107.2675825362452
$ws.Range("R9").Value = 965.4082428262071
$ws.Range("S9").Value = 0.01125638908465337
$ws.Range("T9").Value = 0.01125638908465337

$ws.Range("G10").Value = 82.73961633333333
$ws.Range("H10").Value = 248.218849
$ws.Range("I10").Value = 0.04753008607316088
$ws.Range("J10").Value = 0.04753008607316087
$ws.Range("M10").Value = 0.4248633333333334
$ws.Range("N10").Value = 1.27459
$ws.Range("O10").Value = 0.07761126413398003
$ws.Range("P10").Value = 0.07761126413398005
$ws.Range("Q10").Value = 35.15302919410112
$ws.Range("R10").Value = 316.37726274691
$ws.Range("S10").Value = 0.003688870064534895
$ws.Range("T10").Value = 0.003688870064534895
